# Auto-generated PowerShell COM-interop script
# Applies the row-level data changes described by the diff to sheet1 of the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Cells.Item(4,3).Value = 0
$ws.Cells.Item(4,4).Value = 45845.77758809161
$ws.Cells.Item(4,5).Value = 0
$ws.Cells.Item(4,6).Value = 45845.73388888889

# Row 5
$ws.Cells.Item(5,3).Value = -8
$ws.Cells.Item(5,4).Value = 45845.77758809567
$ws.Cells.Item(5,5).Value = -8
$ws.Cells.Item(5,6).Value = 45845.74221064815

# Row 12
$ws.Cells.Item(12,3).Value = 53
$ws.Cells.Item(12,4).Value = 45845.77758807645
$ws.Cells.Item(12,5).Value = 53
$ws.Cells.Item(12,6).Value = 45845.64517361111

# Row 24
$ws.Cells.Item(24,3).Value = 31
$ws.Cells.Item(24,4).Value = 45845.77758809595
$ws.Cells.Item(24,5).Value = 31
$ws.Cells.Item(24,6).Value = 45845.74221064815

# Row 33
$ws.Cells.Item(33,3).Value = 2499
$ws.Cells.Item(33,4).Value = 45845.77756810876
$ws.Cells.Item(33,5).Value = 2499
$ws.Cells.Item(33,6).Value = 45845.45890046296

# Row 42
$ws.Cells.Item(42,3).Value = -11
$ws.Cells.Item(42,4).Value = 45845.7775681143
$ws.Cells.Item(42,5).Value = -11
$ws.Cells.Item(42,6).Value = 45845.52030092593

# Row 49
$ws.Cells.Item(49,3).Value = 59
$ws.Cells.Item(49,4).Value = 45845.77758809621
$ws.Cells.Item(49,5).Value = 59
$ws.Cells.Item(49,6).Value = 45845.74221064815

# Row 52
$ws.Cells.Item(52,3).Value = 9
$ws.Cells.Item(52,4).Value = 45845.7775880919
$ws.Cells.Item(52,5).Value = 9
$ws.Cells.Item(52,6).Value = 45845.73388888889

# Row 56
$ws.Cells.Item(56,3).Value = 155
$ws.Cells.Item(56,4).Value = 45845.77758809219
$ws.Cells.Item(56,5).Value = 155
$ws.Cells.Item(56,6).Value = 45845.73388888889

# Row 58
$ws.Cells.Item(58,3).Value = 46
$ws.Cells.Item(58,4).Value = 45845.77756810971
$ws.Cells.Item(58,5).Value = 46
$ws.Cells.Item(58,6).Value = 45845.46778935185

# Row 59
$ws.Cells.Item(59,3).Value = 71
$ws.Cells.Item(59,4).Value = 45845.77758808045
$ws.Cells.Item(59,5).Value = 71
$ws.Cells.Item(59,6).Value = 45845.65998842593

# Row 63
$ws.Cells.Item(63,3).Value = 71
$ws.Cells.Item(63,4).Value = 45845.77758809245
$ws.Cells.Item(63,5).Value = 71
$ws.Cells.Item(63,6).Value = 45845.73388888889

# Row 70
$ws.Cells.Item(70,3).Value = 21
$ws.Cells.Item(70,4).Value = 45845.7775681101
$ws.Cells.Item(70,5).Value = 21
$ws.Cells.Item(70,6).Value = 45845.46778935185

# Row 81
$ws.Cells.Item(81,3).Value = 133
$ws.Cells.Item(81,4).Value = 45845.77758809651
$ws.Cells.Item(81,5).Value = 133
$ws.Cells.Item(81,6).Value = 45845.74221064815

# Row 86
$ws.Cells.Item(86,3).Value = 82
$ws.Cells.Item(86,4).Value = 45845.7775880968
$ws.Cells.Item(86,5).Value = 82
$ws.Cells.Item(86,6).Value = 45845.74221064815

# Row 87
$ws.Cells.Item(87,3).Value = -26
$ws.Cells.Item(87,4).Value = 45845.77758807471
$ws.Cells.Item(87,5).Value = -26
$ws.Cells.Item(87,6).Value = 45845.52243055555

# Row 94
$ws.Cells.Item(94,3).Value = 144
$ws.Cells.Item(94,4).Value = 45845.77758809132
$ws.Cells.Item(94,5).Value = 144
$ws.Cells.Item(94,6).Value = 45845.68293981482

# Row 98
$ws.Cells.Item(98,3).Value = 511
$ws.Cells.Item(98,4).Value = 45845.77758811195
$ws.Cells.Item(98,5).Value = 511
$ws.Cells.Item(98,6).Value = 45845.74690972222

# Row 101
$ws.Cells.Item(101,3).Value = 788
$ws.Cells.Item(101,4).Value = 45845.7775880972
$ws.Cells.Item(101,5).Value = 788
$ws.Cells.Item(101,6).Value = 45845.74221064815

# Row 106
$ws.Cells.Item(106,3).Value = 50
$ws.Cells.Item(106,4).Value = 45845.77758808222
$ws.Cells.Item(106,5).Value = 50
$ws.Cells.Item(106,6).Value = 45845.67037037037

# Row 117
$ws.Cells.Item(117,3).Value = 968
$ws.Cells.Item(117,4).Value = 45845.77756811042
$ws.Cells.Item(117,5).Value = 968
$ws.Cells.Item(117,6).Value = 45845.46778935185

# Row 120
$ws.Cells.Item(120,3).Value = 7
$ws.Cells.Item(120,4).Value = 45845.77758808249
$ws.Cells.Item(120,5).Value = 7
$ws.Cells.Item(120,6).Value = 45845.67037037037

# Row 123
$ws.Cells.Item(123,3).Value = 321
$ws.Cells.Item(123,4).Value = 45845.77758807682
$ws.Cells.Item(123,5).Value = 321
$ws.Cells.Item(123,6).Value = 45845.65517361111

# Row 124
$ws.Cells.Item(124,3).Value = 226
$ws.Cells.Item(124,4).Value = 45845.77758809274
$ws.Cells.Item(124,5).Value = 226
$ws.Cells.Item(124,6).Value = 45845.73388888889

# Row 139
$ws.Cells.Item(139,3).Value = 75
$ws.Cells.Item(139,4).Value = 45845.77756810338
$ws.Cells.Item(139,5).Value = 75
$ws.Cells.Item(139,6).Value = 45845.45762731481

# Row 141
$ws.Cells.Item(141,3).Value = 240
$ws.Cells.Item(141,4).Value = 45845.77758810235
$ws.Cells.Item(141,5).Value = 240
$ws.Cells.Item(141,6).Value = 45845.74221064815

# Row 151
$ws.Cells.Item(151,3).Value = -2
$ws.Cells.Item(151,4).Value = 45845.77758807709
$ws.Cells.Item(151,5).Value = -2
$ws.Cells.Item(151,6).Value = 45845.65517361111

# Row 185
$ws.Cells.Item(185,3).Value = 86
$ws.Cells.Item(185,4).Value = 45845.77756810478
$ws.Cells.Item(185,5).Value = 86
$ws.Cells.Item(185,6).Value = 45845.45791666667

# Row 191
$ws.Cells.Item(191,3).Value = 451
$ws.Cells.Item(191,4).Value = 45845.77756810366
$ws.Cells.Item(191,5).Value = 451
$ws.Cells.Item(191,6).Value = 45845.45762731481

# Row 192
$ws.Cells.Item(192,3).Value = 12
$ws.Cells.Item(192,4).Value = 45845.77758810324
$ws.Cells.Item(192,5).Value = 12
$ws.Cells.Item(192,6).Value = 45845.74221064815

# Row 195
$ws.Cells.Item(195,3).Value = -3
$ws.Cells.Item(195,4).Value = 45845.77756808941
$ws.Cells.Item(195,5).Value = -3
$ws.Cells.Item(195,6).Value = 45845.44584490741

# Row 218
$ws.Cells.Item(218,3).Value = 41
$ws.Cells.Item(218,4).Value = 45845.77758810397
$ws.Cells.Item(218,5).Value = 41
$ws.Cells.Item(218,6).Value = 45845.74221064815

# Row 235
$ws.Cells.Item(235,3).Value = 120
$ws.Cells.Item(235,4).Value = 45845.77756809506
$ws.Cells.Item(235,5).Value = 120
$ws.Cells.Item(235,6).Value = 45845.44905092593

# Row 247
$ws.Cells.Item(247,3).Value = 649
$ws.Cells.Item(247,4).Value = 45845.77758809304
$ws.Cells.Item(247,5).Value = 649
$ws.Cells.Item(247,6).Value = 45845.73388888889

# Row 255
$ws.Cells.Item(255,3).Value = 1048
$ws.Cells.Item(255,4).Value = 45845.7775881047
$ws.Cells.Item(255,5).Value = 1048
$ws.Cells.Item(255,6).Value = 45845.74221064815

# Row 258
$ws.Cells.Item(258,3).Value = 55
$ws.Cells.Item(258,4).Value = 45845.77758810546
$ws.Cells.Item(258,5).Value = 55
$ws.Cells.Item(258,6).Value = 45845.74221064815

# Row 272
$ws.Cells.Item(272,3).Value = 194
$ws.Cells.Item(272,4).Value = 45845.77756809169
$ws.Cells.Item(272,5).Value = 194
$ws.Cells.Item(272,6).Value = 45845.44630787037

# Row 273
$ws.Cells.Item(273,3).Value = -2
$ws.Cells.Item(273,4).Value = 45845.77756810197
$ws.Cells.Item(273,5).Value = -2
$ws.Cells.Item(273,6).Value = 45845.45733796297

# Row 274
$ws.Cells.Item(274,3).Value = 198
$ws.Cells.Item(274,4).Value = 45845.77756809197
$ws.Cells.Item(274,5).Value = 198
$ws.Cells.Item(274,6).Value = 45845.44630787037

# Row 283
$ws.Cells.Item(283,3).Value = 114
$ws.Cells.Item(283,4).Value = 45845.77758808278
$ws.Cells.Item(283,5).Value = 114
$ws.Cells.Item(283,6).Value = 45845.67037037037

# Row 287
$ws.Cells.Item(287,3).Value = 414
$ws.Cells.Item(287,4).Value = 45845.77756809338
$ws.Cells.Item(287,5).Value = 414
$ws.Cells.Item(287,6).Value = 45845.44658564815

# Row 291
$ws.Cells.Item(291,3).Value = 405
$ws.Cells.Item(291,4).Value = 45845.77756811069
$ws.Cells.Item(291,5).Value = 405
$ws.Cells.Item(291,6).Value = 45845.46778935185

# Row 295
$ws.Cells.Item(295,3).Value = -33
$ws.Cells.Item(295,4).Value = 45845.77758808305
$ws.Cells.Item(295,5).Value = -33
$ws.Cells.Item(295,6).Value = 45845.67037037037

# Row 309
$ws.Cells.Item(309,3).Value = 904
$ws.Cells.Item(309,4).Value = 45845.77758808334
$ws.Cells.Item(309,5).Value = 904
$ws.Cells.Item(309,6).Value = 45845.67037037037

# Row 314
$ws.Cells.Item(314,3).Value = 13
$ws.Cells.Item(314,4).Value = 45845.77756810081
$ws.Cells.Item(314,5).Value = 13
$ws.Cells.Item(314,6).Value = 45845.45145833334
$ws.Cells.Item(314,7).Value = 0
$ws.Cells.Item(314,8).Value = "Consistente"

# Row 315
$ws.Cells.Item(315,3).Value = 72
$ws.Cells.Item(315,4).Value = 45845.77756808826
$ws.Cells.Item(315,5).Value = 72
$ws.Cells.Item(315,6).Value = 45845.44150462963

# Row 318
$ws.Cells.Item(318,3).Value = 924
$ws.Cells.Item(318,4).Value = 45845.77758810617
$ws.Cells.Item(318,5).Value = 924
$ws.Cells.Item(318,6).Value = 45845.74221064815

# Row 326
$ws.Cells.Item(326,3).Value = 1033
$ws.Cells.Item(326,4).Value = 45845.77758808361
$ws.Cells.Item(326,5).Value = 1033
$ws.Cells.Item(326,6).Value = 45845.67037037037

# Row 328
$ws.Cells.Item(328,3).Value = 171
$ws.Cells.Item(328,4).Value = 45845.77758810679
$ws.Cells.Item(328,5).Value = 171
$ws.Cells.Item(328,6).Value = 45845.74221064815

# Row 342
$ws.Cells.Item(342,3).Value = 50
$ws.Cells.Item(342,4).Value = 45845.77756809365
$ws.Cells.Item(342,5).Value = 50
$ws.Cells.Item(342,6).Value = 45845.44658564815

# Row 346
$ws.Cells.Item(346,3).Value = 12
$ws.Cells.Item(346,4).Value = 45845.77758810715
$ws.Cells.Item(346,5).Value = 12
$ws.Cells.Item(346,6).Value = 45845.74221064815

# Row 350
$ws.Cells.Item(350,3).Value = -30
$ws.Cells.Item(350,4).Value = 45845.77756810395
$ws.Cells.Item(350,5).Value = -30
$ws.Cells.Item(350,6).Value = 45845.45762731481

# Row 354
$ws.Cells.Item(354,3).Value = -4
$ws.Cells.Item(354,4).Value = 45845.77758810748
$ws.Cells.Item(354,5).Value = -4
$ws.Cells.Item(354,6).Value = 45845.74221064815

# Row 358
$ws.Cells.Item(358,3).Value = 21
$ws.Cells.Item(358,4).Value = 45845.77758810783
$ws.Cells.Item(358,5).Value = 21
$ws.Cells.Item(358,6).Value = 45845.74221064815

# Row 361
$ws.Cells.Item(361,3).Value = 363
$ws.Cells.Item(361,4).Value = 45845.77758808082
$ws.Cells.Item(361,5).Value = 363
$ws.Cells.Item(361,6).Value = 45845.65998842593

# Row 363
$ws.Cells.Item(363,3).Value = 313
$ws.Cells.Item(363,4).Value = 45845.77758810822
$ws.Cells.Item(363,5).Value = 313
$ws.Cells.Item(363,6).Value = 45845.74221064815

# Row 371
$ws.Cells.Item(371,3).Value = 3
$ws.Cells.Item(371,4).Value = 45845.77756809533
$ws.Cells.Item(371,5).Value = 3
$ws.Cells.Item(371,6).Value = 45845.44905092593

# Row 385
$ws.Cells.Item(385,3).Value = 208
$ws.Cells.Item(385,4).Value = 45845.77756808967
$ws.Cells.Item(385,5).Value = 208
$ws.Cells.Item(385,6).Value = 45845.44584490741

# Row 390
$ws.Cells.Item(390,3).Value = 137
$ws.Cells.Item(390,4).Value = 45845.7775681017
$ws.Cells.Item(390,5).Value = 137
$ws.Cells.Item(390,6).Value = 45845.45549768519

# Row 394
$ws.Cells.Item(394,3).Value = 148
$ws.Cells.Item(394,4).Value = 45845.77756811484
$ws.Cells.Item(394,5).Value = 148
$ws.Cells.Item(394,6).Value = 45845.52229166667

# Row 395
$ws.Cells.Item(395,3).Value = 35
$ws.Cells.Item(395,4).Value = 45845.77756808996
$ws.Cells.Item(395,5).Value = 35
$ws.Cells.Item(395,6).Value = 45845.44584490741

# Row 402
$ws.Cells.Item(402,3).Value = 78
$ws.Cells.Item(402,4).Value = 45845.77756811097
$ws.Cells.Item(402,5).Value = 78
$ws.Cells.Item(402,6).Value = 45845.46778935185

# Row 403
$ws.Cells.Item(403,3).Value = 37
$ws.Cells.Item(403,4).Value = 45845.77758810855
$ws.Cells.Item(403,5).Value = 37
$ws.Cells.Item(403,6).Value = 45845.74221064815

# Row 404
$ws.Cells.Item(404,3).Value = 419
$ws.Cells.Item(404,4).Value = 45845.77758809332
$ws.Cells.Item(404,5).Value = 419
$ws.Cells.Item(404,6).Value = 45845.73388888889

# Row 410
$ws.Cells.Item(410,3).Value = 1305
$ws.Cells.Item(410,4).Value = 45845.77756810535
$ws.Cells.Item(410,5).Value = 1305
$ws.Cells.Item(410,6).Value = 45845.45826388889

# Row 416
$ws.Cells.Item(416,3).Value = 40
$ws.Cells.Item(416,4).Value = 45845.77756808914
$ws.Cells.Item(416,5).Value = 40
$ws.Cells.Item(416,6).Value = 45845.44221064815

# Row 418
$ws.Cells.Item(418,3).Value = 689
$ws.Cells.Item(418,4).Value = 45845.77758808387
$ws.Cells.Item(418,5).Value = 689
$ws.Cells.Item(418,6).Value = 45845.67037037037

# Row 422
$ws.Cells.Item(422,3).Value = 29
$ws.Cells.Item(422,4).Value = 45845.77756809226
$ws.Cells.Item(422,5).Value = 29
$ws.Cells.Item(422,6).Value = 45845.44630787037

# Row 435
$ws.Cells.Item(435,3).Value = 6
$ws.Cells.Item(435,4).Value = 45845.77756809392
$ws.Cells.Item(435,5).Value = 6
$ws.Cells.Item(435,6).Value = 45845.44658564815

# Row 439
$ws.Cells.Item(439,3).Value = 4
$ws.Cells.Item(439,4).Value = 45845.777588075
$ws.Cells.Item(439,5).Value = 4
$ws.Cells.Item(439,6).Value = 45845.52243055555

# Row 441
$ws.Cells.Item(441,3).Value = 2
$ws.Cells.Item(441,4).Value = 45845.77758811557
$ws.Cells.Item(441,5).Value = 2
$ws.Cells.Item(441,6).Value = 45845.76883101852

# Row 461
$ws.Cells.Item(461,3).Value = 127
$ws.Cells.Item(461,4).Value = 45845.77756810903
$ws.Cells.Item(461,5).Value = 127
$ws.Cells.Item(461,6).Value = 45845.45890046296

# Row 465
$ws.Cells.Item(465,3).Value = 4
$ws.Cells.Item(465,4).Value = 45845.77758811424
$ws.Cells.Item(465,5).Value = 4
$ws.Cells.Item(465,6).Value = 45845.76450231481

# Row 469
$ws.Cells.Item(469,3).Value = 2835
$ws.Cells.Item(469,4).Value = 45845.77756810564
$ws.Cells.Item(469,5).Value = 2835
$ws.Cells.Item(469,6).Value = 45845.45826388889

# Row 480
$ws.Cells.Item(480,3).Value = 199
$ws.Cells.Item(480,4).Value = 45845.77758808112
$ws.Cells.Item(480,5).Value = 199
$ws.Cells.Item(480,6).Value = 45845.65998842593

# Row 485
$ws.Cells.Item(485,3).Value = 6
$ws.Cells.Item(485,4).Value = 45845.77756809024
$ws.Cells.Item(485,5).Value = 6
$ws.Cells.Item(485,6).Value = 45845.44584490741

# Row 488
$ws.Cells.Item(488,4).Value = 45845.7775880897

# Row 489
$ws.Cells.Item(489,4).Value = 45845.77758808924

# Row 490
$ws.Cells.Item(490,4).Value = 45845.77758809012

# Row 493
$ws.Cells.Item(493,3).Value = 523
$ws.Cells.Item(493,4).Value = 45845.77756809563
$ws.Cells.Item(493,5).Value = 523
$ws.Cells.Item(493,6).Value = 45845.44905092593

# Row 496
$ws.Cells.Item(496,3).Value = 12
$ws.Cells.Item(496,4).Value = 45845.7775681151
$ws.Cells.Item(496,5).Value = 12
$ws.Cells.Item(496,6).Value = 45845.52229166667

# Row 507
$ws.Cells.Item(507,3).Value = 55
$ws.Cells.Item(507,4).Value = 45845.77758811493
$ws.Cells.Item(507,5).Value = 55
$ws.Cells.Item(507,6).Value = 45845.76450231481

# Row 510
$ws.Cells.Item(510,3).Value = 132
$ws.Cells.Item(510,4).Value = 45845.77758809361
$ws.Cells.Item(510,5).Value = 132
$ws.Cells.Item(510,6).Value = 45845.73388888889

# Row 527
$ws.Cells.Item(527,3).Value = 10
$ws.Cells.Item(527,4).Value = 45845.77758810887
$ws.Cells.Item(527,5).Value = 10
$ws.Cells.Item(527,6).Value = 45845.74221064815

# Row 535
$ws.Cells.Item(535,3).Value = 84
$ws.Cells.Item(535,4).Value = 45845.77756808857
$ws.Cells.Item(535,5).Value = 84
$ws.Cells.Item(535,6).Value = 45845.44181712963

# Row 542
$ws.Cells.Item(542,3).Value = 75
$ws.Cells.Item(542,4).Value = 45845.77758809538
$ws.Cells.Item(542,5).Value = 75
$ws.Cells.Item(542,6).Value = 45845.52229166667

# Row 545
$ws.Cells.Item(545,3).Value = 4
$ws.Cells.Item(545,4).Value = 45845.77756810593
$ws.Cells.Item(545,5).Value = 4
$ws.Cells.Item(545,6).Value = 45845.45826388889

# Row 569
$ws.Cells.Item(569,3).Value = -2
$ws.Cells.Item(569,4).Value = 45845.77756809421
$ws.Cells.Item(569,5).Value = -2
$ws.Cells.Item(569,6).Value = 45845.44658564815

# Row 570
$ws.Cells.Item(570,3).Value = 2442
$ws.Cells.Item(570,4).Value = 45845.77758807738
$ws.Cells.Item(570,5).Value = 2442
$ws.Cells.Item(570,6).Value = 45845.65517361111

# Row 631
$ws.Cells.Item(631,3).Value = 27
$ws.Cells.Item(631,4).Value = 45845.77756811537
$ws.Cells.Item(631,5).Value = 27
$ws.Cells.Item(631,6).Value = 45845.52229166667

# Row 634
$ws.Cells.Item(634,3).Value = 15
$ws.Cells.Item(634,4).Value = 45845.77756810623
$ws.Cells.Item(634,5).Value = 15
$ws.Cells.Item(634,6).Value = 45845.45826388889

# Row 652
$ws.Cells.Item(652,3).Value = 2
$ws.Cells.Item(652,4).Value = 45845.77756810225
$ws.Cells.Item(652,5).Value = 2
$ws.Cells.Item(652,6).Value = 45845.45733796297

# Row 657
$ws.Cells.Item(657,3).Value = 1715
$ws.Cells.Item(657,4).Value = 45845.7775881092
$ws.Cells.Item(657,5).Value = 1715
$ws.Cells.Item(657,6).Value = 45845.74221064815

# Row 681
$ws.Cells.Item(681,3).Value = -26
$ws.Cells.Item(681,4).Value = 45845.7775681065
$ws.Cells.Item(681,5).Value = -26
$ws.Cells.Item(681,6).Value = 45845.45826388889

# Row 691
$ws.Cells.Item(691,3).Value = 17
$ws.Cells.Item(691,4).Value = 45845.77756811126
$ws.Cells.Item(691,5).Value = 17
$ws.Cells.Item(691,6).Value = 45845.46778935185

# Row 701
$ws.Cells.Item(701,3).Value = 74
$ws.Cells.Item(701,4).Value = 45845.77756809253
$ws.Cells.Item(701,5).Value = 74
$ws.Cells.Item(701,6).Value = 45845.44630787037

# Row 716
$ws.Cells.Item(716,3).Value = 21
$ws.Cells.Item(716,4).Value = 45845.77756811156
$ws.Cells.Item(716,5).Value = 21
$ws.Cells.Item(716,6).Value = 45845.46778935185

# Row 717
$ws.Cells.Item(717,3).Value = 129
$ws.Cells.Item(717,4).Value = 45845.77756809592
$ws.Cells.Item(717,5).Value = 129
$ws.Cells.Item(717,6).Value = 45845.44905092593

# Row 720
$ws.Cells.Item(720,3).Value = 371
$ws.Cells.Item(720,4).Value = 45845.77756811563
$ws.Cells.Item(720,5).Value = 371
$ws.Cells.Item(720,6).Value = 45845.52229166667

# Row 725
$ws.Cells.Item(725,3).Value = 4
$ws.Cells.Item(725,4).Value = 45845.77756810113
$ws.Cells.Item(725,5).Value = 4
$ws.Cells.Item(725,6).Value = 45845.45179398148

# Row 726
$ws.Cells.Item(726,3).Value = 16
$ws.Cells.Item(726,4).Value = 45845.7775880796
$ws.Cells.Item(726,5).Value = 16
$ws.Cells.Item(726,6).Value = 45845.6571412037

# Row 729
$ws.Cells.Item(729,3).Value = 5
$ws.Cells.Item(729,4).Value = 45845.77756810254
$ws.Cells.Item(729,5).Value = 5
$ws.Cells.Item(729,6).Value = 45845.45733796297

# Row 730
$ws.Cells.Item(730,3).Value = 74
$ws.Cells.Item(730,4).Value = 45845.77756809281
$ws.Cells.Item(730,5).Value = 74
$ws.Cells.Item(730,6).Value = 45845.44630787037

# Row 731
$ws.Cells.Item(731,3).Value = 60
$ws.Cells.Item(731,4).Value = 45845.77756810141
$ws.Cells.Item(731,5).Value = 60
$ws.Cells.Item(731,6).Value = 45845.45179398148

# Row 732
$ws.Cells.Item(732,3).Value = 140
$ws.Cells.Item(732,4).Value = 45845.77756810681
$ws.Cells.Item(732,5).Value = 140
$ws.Cells.Item(732,6).Value = 45845.45826388889

# Row 735
$ws.Cells.Item(735,3).Value = 34
$ws.Cells.Item(735,4).Value = 45845.77756810709
$ws.Cells.Item(735,5).Value = 34
$ws.Cells.Item(735,6).Value = 45845.45826388889

# Row 741
$ws.Cells.Item(741,3).Value = -9
$ws.Cells.Item(741,4).Value = 45845.77756810736
$ws.Cells.Item(741,5).Value = -9
$ws.Cells.Item(741,6).Value = 45845.45826388889

# Row 744
$ws.Cells.Item(744,3).Value = 1
$ws.Cells.Item(744,4).Value = 45845.77756809449
$ws.Cells.Item(744,5).Value = 1
$ws.Cells.Item(744,6).Value = 45845.44658564815

# Row 757
$ws.Cells.Item(757,3).Value = 95
$ws.Cells.Item(757,4).Value = 45845.77758809072
$ws.Cells.Item(757,5).Value = 95
$ws.Cells.Item(757,6).Value = 45845.68189814815

# Row 797
$ws.Cells.Item(797,3).Value = -3
$ws.Cells.Item(797,4).Value = 45845.77756810766
$ws.Cells.Item(797,5).Value = -3
$ws.Cells.Item(797,6).Value = 45845.45826388889

# Row 798
$ws.Cells.Item(798,3).Value = 26
$ws.Cells.Item(798,4).Value = 45845.77756810282
$ws.Cells.Item(798,5).Value = 26
$ws.Cells.Item(798,6).Value = 45845.45733796297

# Row 812
$ws.Cells.Item(812,3).Value = -7
$ws.Cells.Item(812,4).Value = 45845.77756810423
$ws.Cells.Item(812,5).Value = -7
$ws.Cells.Item(812,6).Value = 45845.45762731481

# Row 813
$ws.Cells.Item(813,3).Value = 1
$ws.Cells.Item(813,4).Value = 45845.77756810795
$ws.Cells.Item(813,5).Value = 1
$ws.Cells.Item(813,6).Value = 45845.45826388889

# Row 816
$ws.Cells.Item(816,3).Value = 0
$ws.Cells.Item(816,4).Value = 45845.77756810822
$ws.Cells.Item(816,5).Value = 0
$ws.Cells.Item(816,6).Value = 45845.45826388889

# Row 819
$ws.Cells.Item(819,3).Value = 60
$ws.Cells.Item(819,4).Value = 45845.77758810952
$ws.Cells.Item(819,5).Value = 60
$ws.Cells.Item(819,6).Value = 45845.74221064815

# Row 821
$ws.Cells.Item(821,3).Value = 86
$ws.Cells.Item(821,4).Value = 45845.77756811184
$ws.Cells.Item(821,5).Value = 86
$ws.Cells.Item(821,6).Value = 45845.46778935185

# Row 826
$ws.Cells.Item(826,3).Value = 148
$ws.Cells.Item(826,4).Value = 45845.77758810986
$ws.Cells.Item(826,5).Value = 148
$ws.Cells.Item(826,6).Value = 45845.74221064815

# Row 844
$ws.Cells.Item(844,3).Value = 7
$ws.Cells.Item(844,4).Value = 45845.77758807766
$ws.Cells.Item(844,5).Value = 7
$ws.Cells.Item(844,6).Value = 45845.65517361111

# Row 851
$ws.Cells.Item(851,3).Value = 19
$ws.Cells.Item(851,4).Value = 45845.77758809509
$ws.Cells.Item(851,5).Value = 19
$ws.Cells.Item(851,6).Value = 45845.74061342593

# Row 852
$ws.Cells.Item(852,3).Value = 16
$ws.Cells.Item(852,4).Value = 45845.77758808414
$ws.Cells.Item(852,5).Value = 16
$ws.Cells.Item(852,6).Value = 45845.67037037037

# Row 853
$ws.Cells.Item(853,3).Value = 7
$ws.Cells.Item(853,4).Value = 45845.77758807795
$ws.Cells.Item(853,5).Value = 7
$ws.Cells.Item(853,6).Value = 45845.65517361111

# Row 854
$ws.Cells.Item(854,3).Value = 228
$ws.Cells.Item(854,4).Value = 45845.77756809934
$ws.Cells.Item(854,5).Value = 228
$ws.Cells.Item(854,6).Value = 45845.44953703704

# Row 857
$ws.Cells.Item(857,3).Value = 19
$ws.Cells.Item(857,4).Value = 45845.77758807588
$ws.Cells.Item(857,5).Value = 19
$ws.Cells.Item(857,6).Value = 45845.55436342592

# Row 866
$ws.Cells.Item(866,3).Value = -84
$ws.Cells.Item(866,4).Value = 45845.77758811634
$ws.Cells.Item(866,5).Value = -84
$ws.Cells.Item(866,6).Value = 45845.76883101852

# Row 869
$ws.Cells.Item(869,3).Value = 0
$ws.Cells.Item(869,4).Value = 45845.77758807989
$ws.Cells.Item(869,5).Value = 0
$ws.Cells.Item(869,6).Value = 45845.65851851852

# Row 872
$ws.Cells.Item(872,3).Value = 416
$ws.Cells.Item(872,4).Value = 45845.77756811212
$ws.Cells.Item(872,5).Value = 416
$ws.Cells.Item(872,6).Value = 45845.46778935185

# Row 883
$ws.Cells.Item(883,3).Value = 407
$ws.Cells.Item(883,4).Value = 45845.77756810932
$ws.Cells.Item(883,5).Value = 407
$ws.Cells.Item(883,6).Value = 45845.45890046296

# Row 888
$ws.Cells.Item(888,3).Value = 392
$ws.Cells.Item(888,4).Value = 45845.7775680962
$ws.Cells.Item(888,5).Value = 392
$ws.Cells.Item(888,6).Value = 45845.44905092593

# Row 891
$ws.Cells.Item(891,3).Value = 4
$ws.Cells.Item(891,4).Value = 45845.77756809051
$ws.Cells.Item(891,5).Value = 4
$ws.Cells.Item(891,6).Value = 45845.44584490741

# Row 899
$ws.Cells.Item(899,3).Value = 19
$ws.Cells.Item(899,4).Value = 45845.77758811341
$ws.Cells.Item(899,5).Value = 19
$ws.Cells.Item(899,6).Value = 45845.76262731481

# Row 909
$ws.Cells.Item(909,3).Value = 2
$ws.Cells.Item(909,4).Value = 45845.7775680908
$ws.Cells.Item(909,5).Value = 2
$ws.Cells.Item(909,6).Value = 45845.44584490741

# Row 925
$ws.Cells.Item(925,3).Value = 74
$ws.Cells.Item(925,4).Value = 45845.77756811457
$ws.Cells.Item(925,5).Value = 74
$ws.Cells.Item(925,6).Value = 45845.52197916667

# Row 943
$ws.Cells.Item(943,3).Value = 26
$ws.Cells.Item(943,4).Value = 45845.77756811589
$ws.Cells.Item(943,5).Value = 26
$ws.Cells.Item(943,6).Value = 45845.52229166667

# Row 961
$ws.Cells.Item(961,3).Value = 22
$ws.Cells.Item(961,4).Value = 45845.77758807824
$ws.Cells.Item(961,5).Value = 22
$ws.Cells.Item(961,6).Value = 45845.65517361111

# Row 963
$ws.Cells.Item(963,3).Value = 1432
$ws.Cells.Item(963,4).Value = 45845.77756809111
$ws.Cells.Item(963,5).Value = 1432
$ws.Cells.Item(963,6).Value = 45845.44584490741

# Row 977
$ws.Cells.Item(977,3).Value = 89
$ws.Cells.Item(977,4).Value = 45845.77756808885
$ws.Cells.Item(977,5).Value = 89
$ws.Cells.Item(977,6).Value = 45845.44181712963

# Row 1002
$ws.Cells.Item(1002,3).Value = 46
$ws.Cells.Item(1002,4).Value = 45845.77756811238
$ws.Cells.Item(1002,5).Value = 46
$ws.Cells.Item(1002,6).Value = 45845.46778935185

# Row 1017
$ws.Cells.Item(1017,3).Value = 499
$ws.Cells.Item(1017,4).Value = 45845.7775880814
$ws.Cells.Item(1017,5).Value = 499
$ws.Cells.Item(1017,6).Value = 45845.65998842593

# Row 1043
$ws.Cells.Item(1043,3).Value = 26
$ws.Cells.Item(1043,4).Value = 45845.77756811265
$ws.Cells.Item(1043,5).Value = 26
$ws.Cells.Item(1043,6).Value = 45845.46778935185

# Row 1061
$ws.Cells.Item(1061,3).Value = 19
$ws.Cells.Item(1061,4).Value = 45845.77758807292
$ws.Cells.Item(1061,5).Value = 19
$ws.Cells.Item(1061,6).Value = 45845.52229166667

# Row 1062
$ws.Cells.Item(1062,3).Value = 94
$ws.Cells.Item(1062,4).Value = 45845.77758808442
$ws.Cells.Item(1062,5).Value = 94
$ws.Cells.Item(1062,6).Value = 45845.67037037037

# Row 1110
$ws.Cells.Item(1110,3).Value = 95
$ws.Cells.Item(1110,4).Value = 45845.7775680931
$ws.Cells.Item(1110,5).Value = 95
$ws.Cells.Item(1110,6).Value = 45845.44630787037

# Row 1122
$ws.Cells.Item(1122,3).Value = 6
$ws.Cells.Item(1122,4).Value = 45845.77754966282
$ws.Cells.Item(1122,5).Value = 6
$ws.Cells.Item(1122,6).Value = 45845.42458333333

# Row 1126
$ws.Cells.Item(1126,3).Value = 516
$ws.Cells.Item(1126,4).Value = 45845.77758807331
$ws.Cells.Item(1126,5).Value = 516
$ws.Cells.Item(1126,6).Value = 45845.52229166667

# Row 1133
$ws.Cells.Item(1133,3).Value = 27
$ws.Cells.Item(1133,4).Value = 45845.77756809648
$ws.Cells.Item(1133,5).Value = 27
$ws.Cells.Item(1133,6).Value = 45845.44905092593

# Row 1134
$ws.Cells.Item(1134,3).Value = -3
$ws.Cells.Item(1134,4).Value = 45845.77758809419
$ws.Cells.Item(1134,5).Value = -3
$ws.Cells.Item(1134,6).Value = 45845.73724537037

# Row 1147
$ws.Cells.Item(1147,3).Value = 5
$ws.Cells.Item(1147,4).Value = 45845.77756810311
$ws.Cells.Item(1147,5).Value = 5
$ws.Cells.Item(1147,6).Value = 45845.45733796297

# Row 1155
$ws.Cells.Item(1155,3).Value = 6
$ws.Cells.Item(1155,4).Value = 45845.77756811291
$ws.Cells.Item(1155,5).Value = 6
$ws.Cells.Item(1155,6).Value = 45845.46778935185

# Row 1167
$ws.Cells.Item(1167,3).Value = 5
$ws.Cells.Item(1167,4).Value = 45845.7775880945
$ws.Cells.Item(1167,5).Value = 5
$ws.Cells.Item(1167,6).Value = 45845.73811342593

# Row 1174
$ws.Cells.Item(1174,3).Value = 61
$ws.Cells.Item(1174,4).Value = 45845.77756809676
$ws.Cells.Item(1174,5).Value = 61
$ws.Cells.Item(1174,6).Value = 45845.44905092593

# Row 1181
$ws.Cells.Item(1181,3).Value = 44
$ws.Cells.Item(1181,4).Value = 45845.7775880738
$ws.Cells.Item(1181,5).Value = 44
$ws.Cells.Item(1181,6).Value = 45845.52229166667

# Row 1183
$ws.Cells.Item(1183,3).Value = 9
$ws.Cells.Item(1183,4).Value = 45845.77756811399
$ws.Cells.Item(1183,5).Value = 9
$ws.Cells.Item(1183,6).Value = 45845.50509259259

# Row 1187
$ws.Cells.Item(1187,3).Value = 39
$ws.Cells.Item(1187,4).Value = 45845.77758809102
$ws.Cells.Item(1187,5).Value = 39
$ws.Cells.Item(1187,6).Value = 45845.68189814815

# Row 1195
$ws.Cells.Item(1195,3).Value = 7
$ws.Cells.Item(1195,4).Value = 45845.77758807556
$ws.Cells.Item(1195,5).Value = 7
$ws.Cells.Item(1195,6).Value = 45845.54243055556

# Row 1199
$ws.Cells.Item(1199,3).Value = 39
$ws.Cells.Item(1199,4).Value = 45845.7775681045
$ws.Cells.Item(1199,5).Value = 39
$ws.Cells.Item(1199,6).Value = 45845.45762731481

# Row 1223
$ws.Cells.Item(1223,3).Value = 119
$ws.Cells.Item(1223,4).Value = 45845.77758807852
$ws.Cells.Item(1223,5).Value = 119
$ws.Cells.Item(1223,6).Value = 45845.65517361111

# Row 1247
$ws.Cells.Item(1247,3).Value = 6
$ws.Cells.Item(1247,4).Value = 45845.77758808469
$ws.Cells.Item(1247,5).Value = 6
$ws.Cells.Item(1247,6).Value = 45845.67037037037

# Row 1253
$ws.Cells.Item(1253,3).Value = 878
$ws.Cells.Item(1253,4).Value = 45845.77758808497
$ws.Cells.Item(1253,5).Value = 878
$ws.Cells.Item(1253,6).Value = 45845.67037037037

# Row 1255
$ws.Cells.Item(1255,3).Value = 6
$ws.Cells.Item(1255,4).Value = 45845.77756809704
$ws.Cells.Item(1255,5).Value = 6
$ws.Cells.Item(1255,6).Value = 45845.44905092593

# Row 1261
$ws.Cells.Item(1261,3).Value = 204
$ws.Cells.Item(1261,4).Value = 45845.77758811261
$ws.Cells.Item(1261,5).Value = 204
$ws.Cells.Item(1261,6).Value = 45845.74690972222

# Row 1275
$ws.Cells.Item(1275,3).Value = 3
$ws.Cells.Item(1275,4).Value = 45845.77758807412
$ws.Cells.Item(1275,5).Value = 3
$ws.Cells.Item(1275,6).Value = 45845.52229166667

# Row 1332
$ws.Cells.Item(1332,3).Value = 5
$ws.Cells.Item(1332,4).Value = 45845.77758807878
$ws.Cells.Item(1332,5).Value = 5
$ws.Cells.Item(1332,6).Value = 45845.65517361111

# Row 1342
$ws.Cells.Item(1342,3).Value = 855
$ws.Cells.Item(1342,4).Value = 45845.77758808524
$ws.Cells.Item(1342,5).Value = 855
$ws.Cells.Item(1342,6).Value = 45845.67037037037

# Row 1385
$ws.Cells.Item(1385,3).Value = 211
$ws.Cells.Item(1385,4).Value = 45845.77756809962
$ws.Cells.Item(1385,5).Value = 211
$ws.Cells.Item(1385,6).Value = 45845.44953703704

# Row 1388
$ws.Cells.Item(1388,3).Value = 121
$ws.Cells.Item(1388,4).Value = 45845.77758811846
$ws.Cells.Item(1388,5).Value = 121
$ws.Cells.Item(1388,6).Value = 45845.77449074074

# Row 1390
$ws.Cells.Item(1390,3).Value = -3
$ws.Cells.Item(1390,4).Value = 45845.7775880939
$ws.Cells.Item(1390,5).Value = -3
$ws.Cells.Item(1390,6).Value = 45845.73388888889

# Row 1391
$ws.Cells.Item(1391,3).Value = 198
$ws.Cells.Item(1391,4).Value = 45845.7775681132
$ws.Cells.Item(1391,5).Value = 198
$ws.Cells.Item(1391,6).Value = 45845.46778935185

# Row 1406
$ws.Cells.Item(1406,6).Value = 45845.7400462963

# Row 1418
$ws.Cells.Item(1418,3).Value = -1
$ws.Cells.Item(1418,4).Value = 45845.77758807441
$ws.Cells.Item(1418,5).Value = -1
$ws.Cells.Item(1418,6).Value = 45845.52229166667

# Row 1421
$ws.Cells.Item(1421,3).Value = 39
$ws.Cells.Item(1421,4).Value = 45845.77756809731
$ws.Cells.Item(1421,5).Value = 39
$ws.Cells.Item(1421,6).Value = 45845.44905092593

# Row 1446
$ws.Cells.Item(1446,3).Value = -14
$ws.Cells.Item(1446,4).Value = 45845.77758807527
$ws.Cells.Item(1446,5).Value = -14
$ws.Cells.Item(1446,6).Value = 45845.52243055555

# Row 1448
$ws.Cells.Item(1448,3).Value = 48
$ws.Cells.Item(1448,4).Value = 45845.77758808016
$ws.Cells.Item(1448,5).Value = 48
$ws.Cells.Item(1448,6).Value = 45845.65958333333

# Row 1483
$ws.Cells.Item(1483,3).Value = 30
$ws.Cells.Item(1483,4).Value = 45845.77756809758
$ws.Cells.Item(1483,5).Value = 30
$ws.Cells.Item(1483,6).Value = 45845.44905092593

# Row 1504
$ws.Cells.Item(1504,3).Value = 2
$ws.Cells.Item(1504,4).Value = 45845.77760534039
$ws.Cells.Item(1504,5).Value = 2
$ws.Cells.Item(1504,6).Value = 45845.78049768518

# Row 1507
$ws.Cells.Item(1507,3).Value = 159
$ws.Cells.Item(1507,4).Value = 45845.77758807907
$ws.Cells.Item(1507,5).Value = 159
$ws.Cells.Item(1507,6).Value = 45845.65517361111

# Row 1509
$ws.Cells.Item(1509,3).Value = 84
$ws.Cells.Item(1509,4).Value = 45845.77756811346
$ws.Cells.Item(1509,5).Value = 84
$ws.Cells.Item(1509,6).Value = 45845.46778935185

# Row 1594
$ws.Cells.Item(1594,3).Value = 47
$ws.Cells.Item(1594,4).Value = 45845.77758811046
$ws.Cells.Item(1594,5).Value = 47
$ws.Cells.Item(1594,6).Value = 45845.74221064815

# Row 1597
$ws.Cells.Item(1597,3).Value = 3953
$ws.Cells.Item(1597,4).Value = 45845.77758811708
$ws.Cells.Item(1597,5).Value = 3953
$ws.Cells.Item(1597,6).Value = 45845.76883101852

# Row 1600
$ws.Cells.Item(1600,3).Value = 18
$ws.Cells.Item(1600,4).Value = 45845.77756809784
$ws.Cells.Item(1600,5).Value = 18
$ws.Cells.Item(1600,6).Value = 45845.44905092593

# Row 1635
$ws.Cells.Item(1635,3).Value = 1710
$ws.Cells.Item(1635,4).Value = 45845.77756809477
$ws.Cells.Item(1635,5).Value = 1710
$ws.Cells.Item(1635,6).Value = 45845.44658564815

# Row 1753
$ws.Cells.Item(1753,3).Value = -2
$ws.Cells.Item(1753,4).Value = 45845.77756809823
$ws.Cells.Item(1753,5).Value = -2
$ws.Cells.Item(1753,6).Value = 45845.44905092593

# Row 1810
$ws.Cells.Item(1810,3).Value = 10
$ws.Cells.Item(1810,4).Value = 45845.77756808761
$ws.Cells.Item(1810,5).Value = 10
$ws.Cells.Item(1810,6).Value = 45845.43465277777

# Row 1824
$ws.Cells.Item(1824,3).Value = 9
$ws.Cells.Item(1824,4).Value = 45845.7775680999
$ws.Cells.Item(1824,5).Value = 9
$ws.Cells.Item(1824,6).Value = 45845.44953703704

# Row 1830
$ws.Cells.Item(1830,4).Value = 45845.77758809043

# Row 1844
$ws.Cells.Item(1844,3).Value = 367
$ws.Cells.Item(1844,4).Value = 45845.77758808851
$ws.Cells.Item(1844,5).Value = 367
$ws.Cells.Item(1844,6).Value = 45845.67037037037

# Row 1855
$ws.Cells.Item(1855,3).Value = 8
$ws.Cells.Item(1855,4).Value = 45845.77756809851
$ws.Cells.Item(1855,5).Value = 8
$ws.Cells.Item(1855,6).Value = 45845.44905092593

# Row 1903
$ws.Cells.Item(1903,3).Value = 6
$ws.Cells.Item(1903,4).Value = 45845.77756809878
$ws.Cells.Item(1903,5).Value = 6
$ws.Cells.Item(1903,6).Value = 45845.44905092593

# Row 1908
$ws.Cells.Item(1908,3).Value = 0
$ws.Cells.Item(1908,4).Value = 45845.77756809906
$ws.Cells.Item(1908,5).Value = 0
$ws.Cells.Item(1908,6).Value = 45845.44905092593

# Row 1941
$ws.Cells.Item(1941,3).Value = 2
$ws.Cells.Item(1941,4).Value = 45845.77756808795
$ws.Cells.Item(1941,5).Value = 2
$ws.Cells.Item(1941,6).Value = 45845.43465277777

# Row 1980
$ws.Cells.Item(1980,6).Value = 45845.78077546296

# Row 2045
$ws.Cells.Item(2045,3).Value = 200
$ws.Cells.Item(2045,4).Value = 45845.77756810848
$ws.Cells.Item(2045,5).Value = 200
$ws.Cells.Item(2045,6).Value = 45845.45826388889

# Row 2124
$ws.Cells.Item(2124,3).Value = 203
$ws.Cells.Item(2124,4).Value = 45845.77756809141
$ws.Cells.Item(2124,5).Value = 203
$ws.Cells.Item(2124,6).Value = 45845.44584490741

# Row 2131
$ws.Cells.Item(2131,3).Value = -1
$ws.Cells.Item(2131,4).Value = 45845.77756810052
$ws.Cells.Item(2131,5).Value = -1
$ws.Cells.Item(2131,6).Value = 45845.44953703704

# Row 2240
$ws.Cells.Item(2240,3).Value = 229
$ws.Cells.Item(2240,4).Value = 45845.77758811104
$ws.Cells.Item(2240,5).Value = 229
$ws.Cells.Item(2240,6).Value = 45845.74221064815

# Row 2241
$ws.Cells.Item(2241,3).Value = 0
$ws.Cells.Item(2241,4).Value = 45845.77758811147
$ws.Cells.Item(2241,5).Value = 0
$ws.Cells.Item(2241,6).Value = 45845.74221064815

# Row 2278
$ws.Cells.Item(2278,6).Value = 45845.77983796296

# Row 2327
$ws.Cells.Item(2327,3).Value = 2
$ws.Cells.Item(2327,4).Value = 45845.77756808722
$ws.Cells.Item(2327,5).Value = 2
$ws.Cells.Item(2327,6).Value = 45845.42458333333

# Row 2406
$ws.Cells.Item(2406,3).Value = 59
$ws.Cells.Item(2406,4).Value = 45845.77756810503
$ws.Cells.Item(2406,5).Value = 59
$ws.Cells.Item(2406,6).Value = 45845.45791666667

# Row 2459
$ws.Cells.Item(2459,3).Value = 11
$ws.Cells.Item(2459,4).Value = 45845.77756811373
$ws.Cells.Item(2459,5).Value = 11
$ws.Cells.Item(2459,6).Value = 45845.46778935185

# Row 2481
$ws.Cells.Item(2481,3).Value = -39
$ws.Cells.Item(2481,4).Value = 45845.77758807934
$ws.Cells.Item(2481,5).Value = -39
$ws.Cells.Item(2481,6).Value = 45845.65517361111

# Row 2484
$ws.Cells.Item(2484,3).Value = -70
$ws.Cells.Item(2484,4).Value = 45845.77758811792
$ws.Cells.Item(2484,5).Value = -70
$ws.Cells.Item(2484,6).Value = 45845.76883101852

# Row 2603
$ws.Cells.Item(2603,6).Value = 45845.66091435185
$ws.Cells.Item(2603,6).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 2604
$ws.Cells.Item(2604,6).Value = 45845.66216435185
$ws.Cells.Item(2604,6).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row 2605
$ws.Cells.Item(2605,1).Value = 44060400
$ws.Cells.Item(2605,2).Value = 1
$ws.Cells.Item(2605,3).Value = 0
$ws.Cells.Item(2605,4).Value = 45845.77760534002
$ws.Cells.Item(2605,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(2605,5).Value = 0
$ws.Cells.Item(2605,6).Value = 45845.78015046296
$ws.Cells.Item(2605,6).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(2605,7).Value = 0
$ws.Cells.Item(2605,8).Value = "Consistente"
